$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.087.04"
$ws.Range("E2").Value = "  -1.35%  "
$ws.Range("D3").Value = "1.792.88"
$ws.Range("E3").Value = "  -0.44%  "
$c = $ws.Range("D4")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$c = $ws.Range("D5")
$c.Value = "'224.78"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "
$c = $ws.Range("D6")
$c.Value = "'0.550"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  -0.01%  "
$c = $ws.Range("D8")
$c.Value = "'32.47"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("E9").Value = "  -1.53%  "
$c = $ws.Range("D10")
$c.Value = "'0.0707"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "2.051.41"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "1.804.70"
$ws.Range("E13").Value = "  +0.14%  "
$c = $ws.Range("D14")
$c.Value = "'10.82"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.07%  "
$c = $ws.Range("D15")
$c.Value = "'0.624"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.51%  "
$ws.Range("D16").Value = "34.055.23"
$ws.Range("E16").Value = "  -1.47%  "
$ws.Range("E17").Value = "  -2.71%  "
$c = $ws.Range("D18")
$c.Value = "'67.96"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.57%  "
$c = $ws.Range("D19")
$c.Value = "'243.20"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("D20").Value = "0.0₃0783"
$ws.Range("E20").Value = "  -1.60%  "
$c = $ws.Range("D21")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "
$c = $ws.Range("D22")
$c.Value = "'10.67"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.30%  "
$ws.Range("E23").Value = "  -3.78%  "
$ws.Range("E24").Value = "  -2.54%  "
$c = $ws.Range("D25")
$c.Value = "'158.69"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.74%  "
$c = $ws.Range("D26")
$c.Value = "'16.22"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("E31").Value = "  +1.27%  "
$c = $ws.Range("D32")
$c.Value = "'3.66"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.30%  "
$ws.Range("E33").Value = "  -3.07%  "
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("D35").Value = "1.387.48"
$ws.Range("E35").Value = "  -3.00%  "
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("E37").Value = "  -2.07%  "
$c = $ws.Range("D38")
$c.Value = "'0.0184"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D39")
$c.Value = "'2.35"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D40")
$c.Value = "'79.08"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -6.25%  "
$ws.Range("E41").Value = "  -3.41%  "
$ws.Range("E42").Value = "  -4.27%  "
$c = $ws.Range("D43")
$c.Value = "'2.16"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").Value = "0.0₆0140"
$ws.Range("E44").Value = "  +11.00%  "
$c = $ws.Range("D45")
$c.Value = "'0.0493"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("E46").Value = "  -0.93%  "
$c = $ws.Range("D47")
$c.Value = "'107.17"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D48")
$c.Value = "'5.84"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -3.13%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.950.14"
$ws.Range("E49").Value = "  -0.24%  "
$c = $ws.Range("D50")
$c.Value = "'0.998"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.26%  "
$c = $ws.Range("D51")
$c.Value = "'11.98"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.34%  "
